# edit.ps1
# Applies the "Outputs/PtX_demand_LU.xlsx" update: rows 2-31 (years 2030/2040/2050
# fuel-group breakdown) are re-ordered to include new "Fossil Gases" and
# "Fossil Liquids" categories, and the table grows from 31 to 37 rows
# (12 fuel-group rows per year instead of 10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("Hydrogen", 2030, $null, $null, $null, [double]"5.283656111244014e-05", $null, [double]"2.316635748153443e-10", [double]"2.459863587056758e-05", $null, $null)
    ,@("Methanol", 2030, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    ,@("Ammonia", 2030, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    ,@("Synthetic Gases", 2030, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    ,@("Biogenic Gases", 2030, $null, $null, [double]"0.0002070302133370287", [double]"1.616663178346188e-05", $null, $null, [double]"4.303970734896349e-06", $null, $null)
    ,@("Fossil Gases", 2030, $null, $null, $null, [double]"0.0002034442821474712", $null, $null, [double]"1.806983956695509e-05", $null, $null)
    ,@("Synthetic Liquids", 2030, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    ,@("Biogenic Liquids", 2030, $null, $null, $null, [double]"0.0009583835277984215", [double]"5.625438470382552e-06", [double]"0.0006338538975374", [double]"0.0008562968738091", [double]"3.332508913526928e-06", [double]"1.782381275225177e-05")
    ,@("Fossil Liquids", 2030, $null, $null, $null, [double]"0.009451835382458173", [double]"3.961831351454017e-05", [double]"0.0057930828534679", [double]"0.005393546397117", [double]"2.017168346457447e-05", [double]"0.0001274279191674")
    ,@("Biomass [Solid]", 2030, $null, $null, [double]"0.0008830572976137143", $null, $null, $null, $null, $null, $null)
    ,@("Renewable Energy Carrier", 2030, $null, $null, [double]"8.694738566017942e-05", $null, $null, $null, $null, $null, $null)
    ,@("Overall Demand", 2030, $null, $null, [double]"0.001177034896610922", [double]"0.01068266638529997", [double]"4.524375198492272e-05", [double]"0.006426936982668875", [double]"0.006296815717098519", [double]"2.35041923781014e-05", [double]"0.0001452517319196518")
    ,@("Hydrogen", 2040, $null, $null, $null, [double]"0.0002538681704911405", $null, [double]"1.939279995423251e-08", [double]"3.554840186828434e-05", $null, $null)
    ,@("Methanol", 2040, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    ,@("Ammonia", 2040, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    ,@("Synthetic Gases", 2040, $null, $null, $null, [double]"9.296301033643273e-11", $null, $null, [double]"1.094815672177908e-11", $null, $null)
    ,@("Biogenic Gases", 2040, $null, $null, [double]"0.0008090434386974429", [double]"1.988819906160648e-05", $null, $null, [double]"7.533376086321518e-06", $null, $null)
    ,@("Fossil Gases", 2040, $null, $null, $null, [double]"0.0001102345306729826", $null, $null, [double]"1.913760421229235e-05", $null, $null)
    ,@("Synthetic Liquids", 2040, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    ,@("Biogenic Liquids", 2040, $null, $null, $null, [double]"0.0004169401728129474", [double]"9.17039276247531e-06", [double]"0.0007716265368593", [double]"0.0005640913233956999", [double]"4.048457469194271e-06", [double]"2.340641662853338e-05")
    ,@("Fossil Liquids", 2040, $null, $null, $null, [double]"0.002578360594290656", [double]"4.259188969353732e-05", [double]"0.0054617946555969", [double]"0.0024037440514657", [double]"1.79005644140739e-05", [double]"0.0001124321287735")
    ,@("Biomass [Solid]", 2040, $null, $null, [double]"0.0009767952413055328", $null, $null, $null, $null, $null, $null)
    ,@("Renewable Energy Carrier", 2040, $null, $null, [double]"0.0003568693389818933", $null, $null, $null, $null, $null, $null)
    ,@("Overall Demand", 2040, $null, $null, [double]"0.002142708018984869", [double]"0.003379291760292343", [double]"5.176228245601263e-05", [double]"0.006233440585256155", [double]"0.003030054767976455", [double]"2.194902188326818e-05", [double]"0.0001358385454020334")
    ,@("Hydrogen", 2050, $null, $null, $null, [double]"0.0003521444160408327", $null, [double]"3.286923386352513e-08", [double]"5.661741231548529e-05", $null, $null)
    ,@("Methanol", 2050, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    ,@("Ammonia", 2050, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    ,@("Synthetic Gases", 2050, $null, $null, $null, [double]"8.70044268609369e-10", $null, $null, [double]"2.68551980919303e-10", $null, $null)
    ,@("Biogenic Gases", 2050, $null, $null, [double]"0.001955551764325928", [double]"3.406830721744977e-06", $null, $null, [double]"2.173206784143539e-06", $null, $null)
    ,@("Fossil Gases", 2050, $null, $null, $null, [double]"6.883362693042184e-06", $null, $null, [double]"7.190442366028786e-06", $null, $null)
    ,@("Synthetic Liquids", 2050, $null, $null, $null, [double]"2.22160304401638e-12", [double]"3.376959998382159e-13", [double]"2.97286737647877e-11", [double]"9.438632566327578e-12", [double]"2.772717958064315e-14", [double]"3.94304446365687e-12")
    ,@("Biogenic Liquids", 2050, $null, $null, $null, [double]"4.141827709932304e-05", [double]"1.63858818076062e-05", [double]"0.0010106645142242", [double]"0.000145301372805227", [double]"5.206957399277256e-06", [double]"3.178818888705797e-05")
    ,@("Fossil Liquids", 2050, $null, $null, $null, [double]"0.0001449010494921588", [double]"3.841608072198006e-05", [double]"0.0049531446154624", [double]"0.0004293438138343", [double]"1.54079268861911e-05", [double]"9.227971885962604e-05")
    ,@("Biomass [Solid]", 2050, $null, $null, [double]"0.001096317097920804", $null, $null, $null, $null, $null, $null)
    ,@("Renewable Energy Carrier", 2050, $null, $null, [double]"0.0009155415375954781", $null, $null, $null, $null, $null, $null)
    ,@("Overall Demand", 2050, $null, $null, [double]"0.003967410399842211", [double]"0.0005487548083129735", [double]"5.480196286728226e-05", [double]"0.005963842028649137", [double]"0.000640626526095798", [double]"2.061488431319554e-05", [double]"0.0001240679116897285")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
